$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by cloning the "2022-Q3" template
#    sheet (same header/style layout), inserted right before it.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Overwrite the header row (same text, just re-asserting to be safe)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Make sure text-like numeric columns (B..G) keep their original
# "stored as text" representation instead of being auto-coerced to
# numbers by the `.Value =` setter.
$q4.Range("B2:G8").NumberFormat = "@"

# Row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "001173"
$q4.Range("C2").Value = "中欧瑾和灵活配置混合 - A"
$q4.Range("D2").Value = "4.56"
$q4.Range("E2").Value = "89.43"
$q4.Range("F2").Value = "3.37"
$q4.Range("G2").Value = "0.1537"
$q4.Range("H2").Value = 10

# Row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "001105"
$q4.Range("C3").Value = "信澳转型创新股票"
$q4.Range("D3").Value = "3.06"
$q4.Range("E3").Value = "93.59"
$q4.Range("F3").Value = "3.46"
$q4.Range("G3").Value = "0.1059"
$q4.Range("H3").Value = 9

# Row 4
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "001174"
$q4.Range("C4").Value = "中欧瑾和灵活配置混合 - C"
$q4.Range("D4").Value = "2.63"
$q4.Range("E4").Value = "89.43"
$q4.Range("F4").Value = "3.37"
$q4.Range("G4").Value = "0.0886"
$q4.Range("H4").Value = 10

# Row 5
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "000166"
$q4.Range("C5").Value = "中海信息产业精选混合"
$q4.Range("D5").Value = "0.72"
$q4.Range("E5").Value = "83.53"
$q4.Range("F5").Value = "3.54"
$q4.Range("G5").Value = "0.0255"
$q4.Range("H5").Value = 7

# Row 6
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "017288"
$q4.Range("C6").Value = "中欧瑾和灵活配置混合 - E"
$q4.Range("D6").Value = "0.45"
$q4.Range("E6").Value = "89.43"
$q4.Range("F6").Value = "3.37"
$q4.Range("G6").Value = "0.0152"
$q4.Range("H6").Value = 10

# Row 7 (new row beyond the 6-row template -> copy col-A formatting first)
$q4.Range("A2").Copy()
$q4.Range("A7").PasteSpecial(-4122)
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "015608"
$q4.Range("C7").Value = "信澳转型创新股票C"
$q4.Range("D7").Value = "0.33"
$q4.Range("E7").Value = "93.59"
$q4.Range("F7").Value = "3.46"
$q4.Range("G7").Value = "0.0114"
$q4.Range("H7").Value = 9

# Row 8 (new row beyond the 6-row template -> copy col-A formatting first)
$q4.Range("A2").Copy()
$q4.Range("A8").PasteSpecial(-4122)
$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "570007"
$q4.Range("C8").Value = "诺德优选30混合"
$q4.Range("D8").Value = "0.19"
$q4.Range("E8").Value = "87.62"
$q4.Range("F8").Value = "4.86"
$q4.Range("G8").Value = "0.0092"
$q4.Range("H8").Value = 10

# ------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the new 2022-Q4 row at
#    the top of the data (row 2) and push the existing rows down.
# ------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("C2").Value = 7
$tot.Range("D2").Value = 0.41

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2022-Q3"
$tot.Range("C3").Value = 5
$tot.Range("D3").Value = 0.29

$tot.Range("A2").Copy()
$tot.Range("A4").PasteSpecial(-4122)
$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2022-Q2"
$tot.Range("C4").Value = 9
$tot.Range("D4").Value = 0.92

# ------------------------------------------------------------------
# 3) Restore the originally-active tab ("2022-Q2" was the selected
#    sheet before this edit); the sheet Copy() above would otherwise
#    leave the brand-new "2022-Q4" sheet selected instead.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Activate()
